$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.876.24'
$ws.Range("E2").Value = '  +0.77%  '
$ws.Range("D3").Value = '1.627.67'
$ws.Range("E3").Value = '  +1.04%  '
$ws.Range("E4").Value = '  +0.70%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.28'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.86%  '
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +0.69%  '
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0900'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.71%  '
$ws.Range("D12").Value = '1.861.85'
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").Value = '1.637.97'
$ws.Range("E13").Value = '  +1.72%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.562'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.21'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.97%  '
$ws.Range("D16").Value = '29.887.43'
$ws.Range("E16").Value = '  +0.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.82'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.41%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.00'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.16'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.17%  '
$ws.Range("E20").Value = '  -0.44%  '
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.11'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.75'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.57'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.42'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.109'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("E30").Value = '  +1.18%  '
$ws.Range("E31").Value = '  +2.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.59%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.16'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("D34").Value = '1.422.79'
$ws.Range("E34").Value = '  -1.09%  '
$ws.Range("E35").Value = '  +4.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.01'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.06%  '
$ws.Range("E37").Value = '  -4.70%  '
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.556'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.45'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.97%  '
$ws.Range("E42").Value = '  -1.55%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.826'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.38%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.97'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("E45").Value = '  +0.66%  '
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("B47").Value = 'RocketPoolETH'
$ws.Range("C47").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D47").Value = '1.766.78'
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.32'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.92'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.46%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '90.82'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.88%  '
$ws.Range("D51").Value = '0.0₆0111'
$ws.Range("E51").Value = '  +9.10%  '
